# KHL Injuries Master Clubs — refresh scrape snapshot.
# - "snapshot": refresh scraped_at (col K) for all still-present rows (2-31),
#   replace the old kunlun (ШДР) row 32 with the new cska (ЦСК) injuries
#   (Бучельников Дмитрий, Моисеев Данила), and re-append the two kunlun
#   players (Бишофф Джейк, Гроло Жереми) as rows 34-35 with fresh timestamps.
# - "new_injured": log the two newly-injured CSKA players.

$wb = $excel.ActiveWorkbook
$snapshot = $wb.Worksheets.Item("snapshot")
$newInjured = $wb.Worksheets.Item("new_injured")

# --- 1. Refresh scraped_at (column K) for existing rows 2-31 ---------------
$scrapedAt = @{
    2  = "2025-12-21T04:24:41.666043+00:00"
    3  = "2025-12-21T04:24:43.735792+00:00"
    4  = "2025-12-21T04:24:43.735825+00:00"
    5  = "2025-12-21T04:24:43.735843+00:00"
    6  = "2025-12-21T04:24:45.413804+00:00"
    7  = "2025-12-21T04:24:47.513089+00:00"
    8  = "2025-12-21T04:24:49.241586+00:00"
    9  = "2025-12-21T04:24:49.241614+00:00"
    10 = "2025-12-21T04:24:51.386745+00:00"
    11 = "2025-12-21T04:24:55.669842+00:00"
    12 = "2025-12-21T04:24:57.297526+00:00"
    13 = "2025-12-21T04:24:59.452330+00:00"
    14 = "2025-12-21T04:25:03.769811+00:00"
    15 = "2025-12-21T04:25:03.769839+00:00"
    16 = "2025-12-21T04:25:03.769861+00:00"
    17 = "2025-12-21T04:25:03.769878+00:00"
    18 = "2025-12-21T04:25:05.875062+00:00"
    19 = "2025-12-21T04:25:05.875085+00:00"
    20 = "2025-12-21T04:25:05.875093+00:00"
    21 = "2025-12-21T04:25:07.575063+00:00"
    22 = "2025-12-21T04:25:07.575090+00:00"
    23 = "2025-12-21T04:25:07.575107+00:00"
    24 = "2025-12-21T04:25:09.230774+00:00"
    25 = "2025-12-21T04:25:09.230802+00:00"
    26 = "2025-12-21T04:25:11.419113+00:00"
    27 = "2025-12-21T04:25:11.419142+00:00"
    28 = "2025-12-21T04:25:11.419159+00:00"
    29 = "2025-12-21T04:25:13.195726+00:00"
    30 = "2025-12-21T04:25:15.288337+00:00"
    31 = "2025-12-21T04:25:15.288364+00:00"
}

foreach ($r in $scrapedAt.Keys) {
    $snapshot.Cells.Item($r, 11).Value = $scrapedAt[$r]
}

# --- 2. Replace row 32 (was kunlun/Бишофф Джейк) with the two new cska
#        injuries, then re-append the two kunlun players with fresh
#        timestamps as rows 34-35 -------------------------------------------
# "number" (E) and "player_id_khl" (G) are digit-only strings in the source
# data — force text format first so COM's type inference doesn't turn them
# into real numbers (the scraper always writes them as plain text).
function Set-SnapshotRow($Row, $TeamAbbr, $TeamName, $TeamSlug, $PlayerName, $Number, $Position, $PlayerIdKhl, $PlayerUid, $Status, $SourceUrl, $ScrapedAt) {
    $snapshot.Cells.Item($Row, 1).Value = $TeamAbbr
    $snapshot.Cells.Item($Row, 2).Value = $TeamName
    $snapshot.Cells.Item($Row, 3).Value = $TeamSlug
    $snapshot.Cells.Item($Row, 4).Value = $PlayerName
    $snapshot.Cells.Item($Row, 5).NumberFormat = "@"
    $snapshot.Cells.Item($Row, 5).Value = $Number
    $snapshot.Cells.Item($Row, 6).Value = $Position
    $snapshot.Cells.Item($Row, 7).NumberFormat = "@"
    $snapshot.Cells.Item($Row, 7).Value = $PlayerIdKhl
    $snapshot.Cells.Item($Row, 8).Value = $PlayerUid
    $snapshot.Cells.Item($Row, 9).Value = $Status
    $snapshot.Cells.Item($Row, 10).Value = $SourceUrl
    $snapshot.Cells.Item($Row, 11).Value = $ScrapedAt
}

Set-SnapshotRow 32 "ЦСК" "ЦСКА" "cska" "Бучельников Дмитрий" "72" "нападающий" "39102" "1369_ЦСК_бучельниковдмитрий" "injured_active" "https://www.khl.ru/clubs/cska/team/" "2025-12-21T04:25:18.682849+00:00"

Set-SnapshotRow 33 "ЦСК" "ЦСКА" "cska" "Моисеев Данила" "93" "нападающий" "23931" "1369_ЦСК_моисеевданила" "injured_active" "https://www.khl.ru/clubs/cska/team/" "2025-12-21T04:25:18.682865+00:00"

Set-SnapshotRow 34 "ШДР" "Драконы" "kunlun" "Бишофф Джейк" "28" "защитник" "45490" "1369_ШДР_бишоффджейк" "injured_active" "https://www.khl.ru/clubs/kunlun/team/" "2025-12-21T04:25:20.835884+00:00"

Set-SnapshotRow 35 "ШДР" "Драконы" "kunlun" "Гроло Жереми" "75" "защитник" "45343" "1369_ШДР_гроложереми" "injured_active" "https://www.khl.ru/clubs/kunlun/team/" "2025-12-21T04:25:20.835907+00:00"

# --- 3. Log the newly-injured cska players in "new_injured" ----------------
# "changed_day" (G) is a plain yyyy-mm-dd string in the source data — force
# text format first so COM's type inference doesn't turn it into a date.
function Set-NewInjuredRow($Row, $TeamAbbr, $TeamName, $PlayerName, $PlayerUid, $Status, $ChangedAt, $ChangedDay) {
    $newInjured.Cells.Item($Row, 1).Value = $TeamAbbr
    $newInjured.Cells.Item($Row, 2).Value = $TeamName
    $newInjured.Cells.Item($Row, 3).Value = $PlayerName
    $newInjured.Cells.Item($Row, 4).Value = $PlayerUid
    $newInjured.Cells.Item($Row, 5).Value = $Status
    $newInjured.Cells.Item($Row, 6).Value = $ChangedAt
    $newInjured.Cells.Item($Row, 7).NumberFormat = "@"
    $newInjured.Cells.Item($Row, 7).Value = $ChangedDay
}

Set-NewInjuredRow 2 "ЦСК" "ЦСКА" "Бучельников Дмитрий" "1369_ЦСК_бучельниковдмитрий" "INJURED_NEW" "2025-12-21T12:25:21.347447+08:00" "2025-12-21"

Set-NewInjuredRow 3 "ЦСК" "ЦСКА" "Моисеев Данила" "1369_ЦСК_моисеевданила" "INJURED_NEW" "2025-12-21T12:25:21.347447+08:00" "2025-12-21"

Write-Output "snapshot rows: $($snapshot.UsedRange.Rows.Count), new_injured rows: $($newInjured.UsedRange.Rows.Count)"
